# Generate Report for Handoff
# - Priority for the four "2382be94 / 57e828b5 / 937a8967 / e6622f98" rows
#   flips from "low" to "ht" on both the zh-cn and de-de handoff sheets.
# - The zh-cn sheet's "Latest Handoff Datetime" for those same rows moves
#   from 2016-09-06 20:41:18 to 2016-09-06 20:41:45.
# - The "Latest HO Xliff Generate Date" / de-de "Latest Handoff Datetime"
#   (shared text) moves from 2016-09-06 20:41:24 to 2016-09-06 20:41:50,
#   which is reflected on the Overview sheet and the de-de sheet alike.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) for rows 4-7: low -> ht
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# zh-cn Latest Handoff Datetime (H) for rows 4-7
$wsZhCn.Range("H4:H7").Value = "2016-09-06 20:41:45"

# Shared "Latest HO Xliff Generate Date" / de-de Latest Handoff Datetime
$wsOverview.Range("G4:G7").Value = "2016-09-06 20:41:50"
$wsDeDe.Range("H4:H7").Value = "2016-09-06 20:41:50"
